$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated PC1 (column B) and PC2 (column C) values for rows 2-17
# reflecting rerun of the Experimental SVM Hyperplane FE observation script.
$ws.Range("B2").Value = 0.02396012198040794
$ws.Range("C2").Value = 0.06400600273917857
$ws.Range("B3").Value = -0.00867595726186674
$ws.Range("C3").Value = -0.004227289390106689
$ws.Range("B4").Value = -0.001855355037936101
$ws.Range("C4").Value = 0.02083425605938381
$ws.Range("B5").Value = -0.05225237839036743
$ws.Range("C5").Value = -0.1926333790171958
$ws.Range("B6").Value = 0.01537435785120165
$ws.Range("C6").Value = 0.1498228947578887
$ws.Range("B7").Value = 0.3466354321183701
$ws.Range("C7").Value = 0.2625994485723737
$ws.Range("B8").Value = 0.2783585600795661
$ws.Range("C8").Value = 0.2604906281470772
$ws.Range("B9").Value = 0.4224983809048624
$ws.Range("C9").Value = -0.1219803733980831
$ws.Range("B10").Value = 0.7061990224518602
$ws.Range("C10").Value = -0.04217937618053209
$ws.Range("B11").Value = -0.01123451454584329
$ws.Range("C11").Value = 0.05215021348341572
$ws.Range("B12").Value = 0.004850931296653562
$ws.Range("C12").Value = 0.3244709696851096
$ws.Range("B13").Value = 0.04581859851258315
$ws.Range("C13").Value = -0.02632104953641389
$ws.Range("B14").Value = 0.2127827154473627
$ws.Range("C14").Value = 0.08353076346980363
$ws.Range("B15").Value = -0.1394899147064791
$ws.Range("C15").Value = 0.7229484033680846
$ws.Range("B16").Value = 0.2091408775902421
$ws.Range("C16").Value = 0.1179366161259491
$ws.Range("B17").Value = -0.1039040292969665
$ws.Range("C17").Value = 0.3608184411995917
